$wb = $excel.ActiveWorkbook

# Old -> new "Date:" / "Time:" values for each of the 28 backward-elimination
# summary sheets (cell B2 on each sheet holds the full OLS text block).
$oldDate = "Sun, 05 Jan 2020"
$newDate = "Wed, 08 Jan 2020"
$oldTimes = @("21:22:43","21:22:43","21:22:43","21:22:43","21:22:43","21:22:43","21:22:43","21:22:43","21:22:43","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44","21:22:44")
$newTimes = @("19:07:46","19:07:46","19:07:46","19:07:46","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47","19:07:47")

$count = $wb.Worksheets.Count
for ($i = 1; $i -le $count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = [string]$cell.Text
    if ($text -and $text.Length -gt 0) {
        $updated = $text.Replace($oldDate, $newDate).Replace($oldTimes[$i-1], $newTimes[$i-1])
        $cell.Value = $updated
    }
}

# Nudge the saved window position (xWindow) to match the re-saved file.
$excel.ActiveWindow.Left = 690
